$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column F: "United States"
$ws.Range("F1").Value = "United States"

# Update numeric values for columns B-F, rows 2-16 (new poll figures after 85% data collected)
$ws.Range("B2").Value = 10.5823754789272
$ws.Range("C2").Value = 14.3776978417266
$ws.Range("D2").Value = 33.2039800995025
$ws.Range("E2").Value = 30.8287937743191
$ws.Range("F2").Value = 14.3871681415929

$ws.Range("B3").Value = 24.1853281853282
$ws.Range("C3").Value = 20.1114864864865
$ws.Range("D3").Value = 21.4766355140187
$ws.Range("E3").Value = 23.3176895306859
$ws.Range("F3").Value = 23.1936542669584

$ws.Range("B4").Value = 13.9748953974895
$ws.Range("C4").Value = 18.2904290429043
$ws.Range("D4").Value = 16.3183856502242
$ws.Range("E4").Value = 10.1853281853282
$ws.Range("F4").Value = 21.7176981541802

$ws.Range("B5").Value = 27.1274131274131
$ws.Range("C5").Value = 16.214983713355
$ws.Range("D5").Value = 25.6869158878505
$ws.Range("E5").Value = 22.4856115107914
$ws.Range("F5").Value = 27.7279821627648

$ws.Range("B6").Value = 12.6984126984127
$ws.Range("C6").Value = 16.5833333333333
$ws.Range("D6").Value = 12.3405172413793
$ws.Range("E6").Value = 23.6339622641509
$ws.Range("F6").Value = 10.0673076923077

$ws.Range("B7").Value = 13.848623853211
$ws.Range("C7").Value = 8.44055944055944
$ws.Range("D7").Value = 9.73684210526316
$ws.Range("E7").Value = 5.39114391143911
$ws.Range("F7").Value = 13.4772486772487

$ws.Range("B8").Value = 11.2057613168724
$ws.Range("C8").Value = 17.5180327868852
$ws.Range("D8").Value = 20.521327014218
$ws.Range("E8").Value = 12.8588709677419
$ws.Range("F8").Value = 13.0299896587384

$ws.Range("B9").Value = 20.8222222222222
$ws.Range("C9").Value = 18.5631067961165
$ws.Range("D9").Value = 14.4691943127962
$ws.Range("E9").Value = 17.972
$ws.Range("F9").Value = 19.6789989118607

$ws.Range("B10").Value = 8.50387596899225
$ws.Range("C10").Value = 9.18210862619808
$ws.Range("D10").Value = 8.87614678899082
$ws.Range("E10").Value = 11.7480916030534
$ws.Range("F10").Value = 10.344502617801

$ws.Range("B11").Value = 15.8143459915612
$ws.Range("C11").Value = 15.5368098159509
$ws.Range("D11").Value = 14.4887892376682
$ws.Range("E11").Value = 15.8978723404255
$ws.Range("F11").Value = 13.2456331877729

$ws.Range("B12").Value = 21.5924369747899
$ws.Range("C12").Value = 17.877133105802
$ws.Range("D12").Value = 16.7920792079208
$ws.Range("E12").Value = 20.1991701244813
$ws.Range("F12").Value = 19.5309917355372

$ws.Range("B13").Value = 20.3899082568807
$ws.Range("C13").Value = 23.9665551839465
$ws.Range("D13").Value = 16.725321888412
$ws.Range("E13").Value = 16.2
$ws.Range("F13").Value = 15.4519337016575

$ws.Range("B14").Value = 20.7725490196078
$ws.Range("C14").Value = 22.5249169435216
$ws.Range("D14").Value = 18.8325123152709
$ws.Range("E14").Value = 19.4961538461538
$ws.Range("F14").Value = 20.1214750542299

$ws.Range("B15").Value = 15.4666666666667
$ws.Range("C15").Value = 16.8954703832753
$ws.Range("D15").Value = 13.5321888412017
$ws.Range("E15").Value = 12.992337164751
$ws.Range("F15").Value = 14.8700102354145

$ws.Range("B16").Value = 12.7015503875969
$ws.Range("C16").Value = 13.9428571428571
$ws.Range("D16").Value = 9.44588744588745
$ws.Range("E16").Value = 5.91836734693878
$ws.Range("F16").Value = 9.01342281879195

